$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert a new blank row at row 38 (shifts old rows 38-41 down to 39-42)
$ws.Rows(38).Insert()

# Step 2: apply all cell value changes to match the target state
$ws.Range('A8').Value = ''
$ws.Range('C8').Value = 93
$ws.Range('D8').Value = '1.0'
$ws.Range('E8').Value = 'Rewiring of light point/ fan point/ exhaust fan point/ call bell point with 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade (IS:694) in recessed  ISI marked MMS ( IS:9537 P - III ) virgin material PVC conduit & it''s ISI marked (IS:3419-1988) accessories, round tiles, 1.2 mm thick MS box with earth terminal, 6 A switch, 3 pin ceiling rose/holder / 3 way connector , 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/  brass  screws, cup washers, making connections, testing etc. as required. For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range('F8').Value = 0
$ws.Range('G8').Value = '0.00'
$ws.Range('C9').Value = 4
$ws.Range('D9').Value = '2'
$ws.Range('E9').Value = 'Short point (up to 3 mtr.)'
$ws.Range('F9').Value = 256
$ws.Range('G9').Value = '1024.00'
$ws.Range('A10').Value = 'P. point'
$ws.Range('C10').Value = 15
$ws.Range('D10').Value = '4'
$ws.Range('E10').Value = 'Long point  (up to 10 mtr.)'
$ws.Range('F10').Value = 662
$ws.Range('G10').Value = '9930.00'
$ws.Range('C11').Value = 25
$ws.Range('G11').Value = '3400.00'
$ws.Range('C12').Value = 56
$ws.Range('D12').Value = '4.0'
$ws.Range('E12').Value = 'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range('F12').Value = 50
$ws.Range('G12').Value = '2800.00'
$ws.Range('C13').Value = 33
$ws.Range('D13').Value = '6.0'
$ws.Range('E13').Value = 'Providing & Fixing of  3/6 pin 16 amp flush type non modular socket  made out from Industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range('F13').Value = 78
$ws.Range('G13').Value = '2574.00'
$ws.Range('D14').Value = '7.0'
$ws.Range('E14').Value = 'Providing & Fixing of  ISI marked (IS:371) 6 amp surface type 3 pin ceiling rose with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screws including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range('F14').Value = 30
$ws.Range('G14').Value = '690.00'
$ws.Range('C15').Value = 49
$ws.Range('D15').Value = '8.0'
$ws.Range('E15').Value = 'Providing & Fixing of ISI marked (IS:1258) batten/angle lamp  holder with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screwsincluding making connection testing etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range('F15').Value = 30
$ws.Range('G15').Value = '1470.00'
$ws.Range('C16').Value = 86
$ws.Range('D16').Value = '9.0'
$ws.Range('E16').Value = 'Providing & Fixing of IS 11037:1984  marked  non modular socket size flush type 180 watt rotary minimum 5 step fan regulator with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range('F16').Value = 219
$ws.Range('G16').Value = '18834.00'
$ws.Range('C17').Value = 66
$ws.Range('D17').Value = '10.0'
$ws.Range('E17').Value = 'Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range('F17').Value = 303
$ws.Range('G17').Value = '19998.00'
$ws.Range('A18').Value = ''
$ws.Range('C18').Value = 70
$ws.Range('D18').Value = '11.0'
$ws.Range('E18').Value = 'S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range('F18').Value = 0
$ws.Range('G18').Value = '0.00'
$ws.Range('A19').Value = 'R. mtr.'
$ws.Range('C19').Value = 37
$ws.Range('D19').Value = '17'
$ws.Range('E19').Value = '25 mm'
$ws.Range('F19').Value = 56
$ws.Range('G19').Value = '2072.00'
$ws.Range('A20').Value = ''
$ws.Range('C20').Value = 66
$ws.Range('D20').Value = '12.0'
$ws.Range('E20').Value = 'Supplying and drawing FR PVC insulated & unsheathed flexible copper conductor as per PWD specification for electrical Works with ISI marked (IS:694) and as per IS 8130 : 2013 of 1.1 kV grade . Wire should be made from  99.90 % purity copper, class 2 stranding in acc. to IS:8130/IEC 60228 for  lower watt loss , oxygen free for less chances of oxidization, insulation PVC type A/C/D , flame retardant as per IS 10810-53, better amperage rating as per IS:3961 part 5,  in existing  surface or recessed PVC/ MS conduit/casing capping making connections with Copper Lugs of suitable size, Ferrules,testing etc. as required. OEM Must have its own in house NABL lab setup for all testing facilities for wires.   For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range('F20').Value = 0
$ws.Range('G20').Value = '0.00'
$ws.Range('A21').Value = 'Mtr.'
$ws.Range('C21').Value = 88
$ws.Range('D21').Value = '19'
$ws.Range('E21').Value = '2 x 2.5 sq. mm. + 1x1.5sqmm'
$ws.Range('F21').Value = 81
$ws.Range('G21').Value = '7128.00'
$ws.Range('A22').Value = 'Mtr.'
$ws.Range('C22').Value = 7
$ws.Range('D22').Value = '20'
$ws.Range('E22').Value = '2 x 4.0 sq. mm. + 1 x 2.5 sq. mm.'
$ws.Range('F22').Value = 122
$ws.Range('G22').Value = '854.00'
$ws.Range('A23').Value = 'Set'
$ws.Range('C23').Value = 57
$ws.Range('D23').Value = '13.0'
$ws.Range('E23').Value = 'Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. ''B'' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range('F23').Value = 5733
$ws.Range('G23').Value = '326781.00'
$ws.Range('A24').Value = ''
$ws.Range('C24').Value = 70
$ws.Range('D24').Value = '14.0'
$ws.Range('E24').Value = 'Supply & Laying following size earth wire in horizontal or vertical run in ground/surface/recess including riveting, soldering, saddles,  making connection with GI/Cu purity purity >95%  thimble etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range('F24').Value = 0
$ws.Range('G24').Value = '0.00'
$ws.Range('A25').Value = 'Mtr.'
$ws.Range('C25').Value = 90
$ws.Range('D25').Value = '23'
$ws.Range('E25').Value = '8 SWG G.I. ( Hot Dipped  ) Wire '
$ws.Range('F25').Value = 20
$ws.Range('G25').Value = '1800.00'
$ws.Range('C26').Value = 51
$ws.Range('A27').Value = ''
$ws.Range('C27').Value = 39
$ws.Range('D27').Value = '16.0'
$ws.Range('E27').Value = 'Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range('F27').Value = 0
$ws.Range('G27').Value = '0.00'
$ws.Range('C28').Value = 38
$ws.Range('D28').Value = '17.0'
$ws.Range('E28').Value = 'Providing & Fixing of 240/415 V AC MCB with positive isolation of 10 kA breaking capacity (B/ C/D tripping characteristic as per type of load and  site requirement) 4 KV impulse withstand voltage, ISI marked IS 8828(1996) / conforming to IEC 60898-1 2002, IEC 60947-2, low watt losses, trip free mechanisum , energy limiting of  class 3 as per IEC,  minimum phase termination capacity of 35sq.mm. , conductor line load reversibility , IP 20 contact protection and fitted in  existing distribution board/sheets, minimum electrical operation 20,000 upto 20 A rating and 10,000 upto 63 A, 5000 for 80 A & above rating  including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range('A29').Value = ''
$ws.Range('C29').Value = 62
$ws.Range('D29').Value = '29'
$ws.Range('E29').Value = 'Single pole MCB   (With B/C curve tripping Characteristics)'
$ws.Range('F29').Value = 0
$ws.Range('G29').Value = '0.00'
$ws.Range('A30').Value = 'Each'
$ws.Range('C30').Value = 13
$ws.Range('D30').Value = '30'
$ws.Range('E30').Value = ' 6 A to 32 A rating'
$ws.Range('F30').Value = 187
$ws.Range('G30').Value = '2431.00'
$ws.Range('C31').Value = 47
$ws.Range('C32').Value = 60
$ws.Range('G32').Value = '54000.00'
$ws.Range('C33').Value = 55
$ws.Range('D33').Value = '18.0'
$ws.Range('E33').Value = 'Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range('A34').Value = ''
$ws.Range('C34').Value = 56
$ws.Range('D34').Value = '34'
$ws.Range('E34').Value = 'Metal door (single phase) IK-09 and IP-43 with Metal end box'
$ws.Range('F34').Value = 0
$ws.Range('G34').Value = '0.00'
$ws.Range('A35').Value = 'Each'
$ws.Range('C35').Value = 29
$ws.Range('D35').Value = '35'
$ws.Range('E35').Value = '8 Way (8+2)'
$ws.Range('F35').Value = 2184
$ws.Range('G35').Value = '63336.00'
$ws.Range('A36').Value = ''
$ws.Range('C36').Value = 10
$ws.Range('D36').Value = '36'
$ws.Range('E36').Value = 'Total'
$ws.Range('A37').Value = '%'
$ws.Range('C37').Value = 78
$ws.Range('D37').Value = '37'
$ws.Range('E37').Value = 'Add Tender Premium '
$ws.Range('G40').Value = '519122.00'
$ws.Range('H40').Value = '519122.00'
$ws.Range('G42').Value = '519122.00'
$ws.Range('H42').Value = '519122.00'
